$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The merged cell B2:B4 held the test case title "Verify Home page";
# rename it to "Verify Landing page".
$ws.Range("B2").Value = "Verify Landing page"

# Update the active selection to match the authored state (B2:B4 selected,
# active cell B2) instead of the stray G15 selection left over from editing.
$ws.Range("B2:B4").Select() | Out-Null
